{"js": "// Replace each three-digit x one-digit multiplication expression\n// with its new value, matching the exact original text (unique per cell).\nconst replacements = [\n  [\"513\u00d78=\", \"508\u00d74=\"],\n  [\"541\u00d73=\", \"398\u00d79=\"],\n  [\"308\u00d72=\", \"967\u00d74=\"],\n  [\"765\u00d79=\", \"225\u00d72=\"],\n  [\"539\u00d72=\", \"219\u00d72=\"],\n  [\"686\u00d74=\", \"545\u00d73=\"],\n  [\"430\u00d74=\", \"794\u00d74=\"],\n  [\"696\u00d79=\", \"317\u00d76=\"],\n  [\"229\u00d74=\", \"360\u00d75=\"],\n  [\"674\u00d77=\", \"993\u00d79=\"],\n  [\"397\u00d76=\", \"442\u00d79=\"],\n  [\"913\u00d78=\", \"257\u00d75=\"],\n  [\"453\u00d76=\", \"133\u00d72=\"],\n  [\"358\u00d72=\", \"322\u00d75=\"],\n  [\"426\u00d73=\", \"728\u00d76=\"],\n  [\"738\u00d77=\", \"846\u00d72=\"],\n  [\"712\u00d76=\", \"162\u00d73=\"],\n  [\"813\u00d75=\", \"682\u00d73=\"],\n  [\"391\u00d75=\", \"920\u00d76=\"],\n  [\"466\u00d76=\", \"992\u00d74=\"],\n  [\"961\u00d73=\", \"621\u00d76=\"],\n  [\"929\u00d74=\", \"589\u00d76=\"],\n  [\"206\u00d78=\", \"134\u00d72=\"],\n  [\"186\u00d78=\", \"962\u00d75=\"],\n  [\"766\u00d77=\", \"321\u00d77=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit x one-digit multiplication expression\n# with its new value. Word Find/Replace, scoped to the whole document\n# body, case-sensitive exact match (each source string is unique).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"513\u00d78=\", \"508\u00d74=\"),\n    @(\"541\u00d73=\", \"398\u00d79=\"),\n    @(\"308\u00d72=\", \"967\u00d74=\"),\n    @(\"765\u00d79=\", \"225\u00d72=\"),\n    @(\"539\u00d72=\", \"219\u00d72=\"),\n    @(\"686\u00d74=\", \"545\u00d73=\"),\n    @(\"430\u00d74=\", \"794\u00d74=\"),\n    @(\"696\u00d79=\", \"317\u00d76=\"),\n    @(\"229\u00d74=\", \"360\u00d75=\"),\n    @(\"674\u00d77=\", \"993\u00d79=\"),\n    @(\"397\u00d76=\", \"442\u00d79=\"),\n    @(\"913\u00d78=\", \"257\u00d75=\"),\n    @(\"453\u00d76=\", \"133\u00d72=\"),\n    @(\"358\u00d72=\", \"322\u00d75=\"),\n    @(\"426\u00d73=\", \"728\u00d76=\"),\n    @(\"738\u00d77=\", \"846\u00d72=\"),\n    @(\"712\u00d76=\", \"162\u00d73=\"),\n    @(\"813\u00d75=\", \"682\u00d73=\"),\n    @(\"391\u00d75=\", \"920\u00d76=\"),\n    @(\"466\u00d76=\", \"992\u00d74=\"),\n    @(\"961\u00d73=\", \"621\u00d76=\"),\n    @(\"929\u00d74=\", \"589\u00d76=\"),\n    @(\"206\u00d78=\", \"134\u00d72=\"),\n    @(\"186\u00d78=\", \"962\u00d75=\"),\n    @(\"766\u00d77=\", \"321\u00d77=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
